$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 3 holds a placeholder "Max Mustermann" test user whose fields were
# still empty/dummy. Fill in the real test values now that the DB hookup
# works end to end.

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "max_mustermann"
$ws.Range("C3").Value = "Max"
$ws.Range("D3").Value = "Mustermann"
$ws.Range("E3").Value = "Musterstraße 1"
$ws.Range("F3").Value = 12345
$ws.Range("G3").Value = "Musterstadt"
$ws.Range("H3").Value = "max.mustermann@example.com"
$ws.Range("I3").Value = "passwort123"
